$wb = $excel.ActiveWorkbook

# Sheet references
$wsPersonnel = $wb.Worksheets.Item(1)   # 인원변동
$wsTurnover  = $wb.Worksheets.Item(2)   # 퇴사율

# Fill in the new zero-value cells on the 퇴사율 sheet
$wsTurnover.Range("D2").Value = 0
$wsTurnover.Range("F2").Value = 0
$wsTurnover.Range("H2").Value = 0
$wsTurnover.Range("G3").Value = 0
$wsTurnover.Range("H3").Value = 0
$wsTurnover.Range("E4").Value = 0

# Move the selection on the first sheet (인원변동) and make it no longer the
# tab-selected sheet
$wsPersonnel.Range("D16").Select() | Out-Null

# Move the selection on the 퇴사율 sheet and make it the active / selected tab
$wsTurnover.Select() | Out-Null
$wsTurnover.Range("E5").Select() | Out-Null

# Make sure the workbook's active sheet (activeTab) is 퇴사율 (index 1)
$wb.Worksheets.Item(2).Activate() | Out-Null
